$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6104469299316406
$ws.Range("B1").Value = 2.796266317367554
$ws.Range("C1").Value = 6.363174438476562
$ws.Range("D1").Value = 1.78946053981781
$ws.Range("E1").Value = 1.572123527526855
